# =====================================================================
# Horarios Linea 141 - actualizacion de scrape (12:44:05 -> 12:57:33)
# Actualiza encabezados, reordena filas con el mismo Hora_Llegada segun
# el nuevo orden de scrape y agrega las filas nuevas al final de cada hoja.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---- Hoja "LP1912": cabecera + filas reordenadas/agregadas (6 -> 264 filas) ----
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2,1).Value = "Última actualización: 12:57:33"
$ws1.Cells.Item(3,1).Value = "Total filas: 264"
$ws1.Cells.Item(54,1).Value = "06:02:16"
$ws1.Cells.Item(54,2).Value = "07:16"
$ws1.Cells.Item(54,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(54,4).Value = 74
$ws1.Cells.Item(54,5).Value = "LP1912"
$ws1.Cells.Item(55,1).Value = "06:37:24"
$ws1.Cells.Item(55,2).Value = "07:16"
$ws1.Cells.Item(55,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(55,4).Value = 39
$ws1.Cells.Item(55,5).Value = "LP1912"
$ws1.Cells.Item(120,1).Value = "08:47:19"
$ws1.Cells.Item(120,2).Value = "09:35"
$ws1.Cells.Item(120,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(120,4).Value = 48
$ws1.Cells.Item(120,5).Value = "LP1912"
$ws1.Cells.Item(121,1).Value = "08:47:19"
$ws1.Cells.Item(121,2).Value = "09:35"
$ws1.Cells.Item(121,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(121,4).Value = 48
$ws1.Cells.Item(121,5).Value = "LP1912"
$ws1.Cells.Item(203,1).Value = "11:15:53"
$ws1.Cells.Item(203,2).Value = "12:20"
$ws1.Cells.Item(203,3).Value = "14_ABASTO"
$ws1.Cells.Item(203,4).Value = 65
$ws1.Cells.Item(203,5).Value = "LP1912"
$ws1.Cells.Item(204,1).Value = "11:15:53"
$ws1.Cells.Item(204,2).Value = "12:20"
$ws1.Cells.Item(204,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(204,4).Value = 65
$ws1.Cells.Item(204,5).Value = "LP1912"
$ws1.Cells.Item(205,1).Value = "10:50:37"
$ws1.Cells.Item(205,2).Value = "12:20"
$ws1.Cells.Item(205,3).Value = "215A_EL PATO"
$ws1.Cells.Item(205,4).Value = 90
$ws1.Cells.Item(205,5).Value = "LP1912"
$ws1.Cells.Item(214,1).Value = "12:24:14"
$ws1.Cells.Item(214,2).Value = "12:35"
$ws1.Cells.Item(214,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(214,4).Value = 11
$ws1.Cells.Item(214,5).Value = "LP1912"
$ws1.Cells.Item(215,1).Value = "12:24:14"
$ws1.Cells.Item(215,2).Value = "12:35"
$ws1.Cells.Item(215,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(215,4).Value = 11
$ws1.Cells.Item(215,5).Value = "LP1912"
$ws1.Cells.Item(223,1).Value = "12:44:05"
$ws1.Cells.Item(223,2).Value = "12:45"
$ws1.Cells.Item(223,3).Value = "10_OLMOS"
$ws1.Cells.Item(223,4).Value = 1
$ws1.Cells.Item(223,5).Value = "LP1912"
$ws1.Cells.Item(224,1).Value = "12:44:05"
$ws1.Cells.Item(224,2).Value = "12:45"
$ws1.Cells.Item(224,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(224,4).Value = 1
$ws1.Cells.Item(224,5).Value = "LP1912"
$ws1.Cells.Item(229,1).Value = "12:57:33"
$ws1.Cells.Item(229,2).Value = "12:57"
$ws1.Cells.Item(229,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(229,4).Value = 0
$ws1.Cells.Item(229,5).Value = "LP1912"
$ws1.Cells.Item(230,1).Value = "12:57:33"
$ws1.Cells.Item(230,2).Value = "12:57"
$ws1.Cells.Item(230,3).Value = "17_ROMERO"
$ws1.Cells.Item(230,4).Value = 0
$ws1.Cells.Item(230,5).Value = "LP1912"
$ws1.Cells.Item(231,1).Value = "12:57:33"
$ws1.Cells.Item(231,2).Value = "12:58"
$ws1.Cells.Item(231,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(231,4).Value = 1
$ws1.Cells.Item(231,5).Value = "LP1912"
$ws1.Cells.Item(232,1).Value = "11:15:53"
$ws1.Cells.Item(232,2).Value = "13:02"
$ws1.Cells.Item(232,3).Value = "15_ABASTO"
$ws1.Cells.Item(232,4).Value = 107
$ws1.Cells.Item(232,5).Value = "LP1912"
$ws1.Cells.Item(233,1).Value = "12:24:14"
$ws1.Cells.Item(233,2).Value = "13:03"
$ws1.Cells.Item(233,3).Value = "14_ABASTO"
$ws1.Cells.Item(233,4).Value = 39
$ws1.Cells.Item(233,5).Value = "LP1912"
$ws1.Cells.Item(234,1).Value = "11:15:53"
$ws1.Cells.Item(234,2).Value = "13:06"
$ws1.Cells.Item(234,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(234,4).Value = 111
$ws1.Cells.Item(234,5).Value = "LP1912"
$ws1.Cells.Item(235,1).Value = "12:57:33"
$ws1.Cells.Item(235,2).Value = "13:06"
$ws1.Cells.Item(235,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(235,4).Value = 9
$ws1.Cells.Item(235,5).Value = "LP1912"
$ws1.Cells.Item(236,1).Value = "12:24:14"
$ws1.Cells.Item(236,2).Value = "13:07"
$ws1.Cells.Item(236,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(236,4).Value = 43
$ws1.Cells.Item(236,5).Value = "LP1912"
$ws1.Cells.Item(237,1).Value = "11:15:53"
$ws1.Cells.Item(237,2).Value = "13:13"
$ws1.Cells.Item(237,3).Value = "215D_EL PATO"
$ws1.Cells.Item(237,4).Value = 118
$ws1.Cells.Item(237,5).Value = "LP1912"
$ws1.Cells.Item(238,1).Value = "11:43:19"
$ws1.Cells.Item(238,2).Value = "13:14"
$ws1.Cells.Item(238,3).Value = "215D_EL PATO"
$ws1.Cells.Item(238,4).Value = 91
$ws1.Cells.Item(238,5).Value = "LP1912"
$ws1.Cells.Item(239,1).Value = "11:43:19"
$ws1.Cells.Item(239,2).Value = "13:14"
$ws1.Cells.Item(239,3).Value = "17_ROMERO"
$ws1.Cells.Item(239,4).Value = 91
$ws1.Cells.Item(239,5).Value = "LP1912"
$ws1.Cells.Item(240,1).Value = "12:57:33"
$ws1.Cells.Item(240,2).Value = "13:14"
$ws1.Cells.Item(240,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(240,4).Value = 17
$ws1.Cells.Item(240,5).Value = "LP1912"
$ws1.Cells.Item(241,1).Value = "12:44:05"
$ws1.Cells.Item(241,2).Value = "13:15"
$ws1.Cells.Item(241,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(241,4).Value = 31
$ws1.Cells.Item(241,5).Value = "LP1912"
$ws1.Cells.Item(242,1).Value = "11:43:19"
$ws1.Cells.Item(242,2).Value = "13:19"
$ws1.Cells.Item(242,3).Value = "10_OLMOS"
$ws1.Cells.Item(242,4).Value = 96
$ws1.Cells.Item(242,5).Value = "LP1912"
$ws1.Cells.Item(243,1).Value = "12:44:05"
$ws1.Cells.Item(243,2).Value = "13:20"
$ws1.Cells.Item(243,3).Value = "10_OLMOS"
$ws1.Cells.Item(243,4).Value = 36
$ws1.Cells.Item(243,5).Value = "LP1912"
$ws1.Cells.Item(244,1).Value = "11:43:19"
$ws1.Cells.Item(244,2).Value = "13:21"
$ws1.Cells.Item(244,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(244,4).Value = 98
$ws1.Cells.Item(244,5).Value = "LP1912"
$ws1.Cells.Item(245,1).Value = "11:43:19"
$ws1.Cells.Item(245,2).Value = "13:26"
$ws1.Cells.Item(245,3).Value = "14_ABASTO"
$ws1.Cells.Item(245,4).Value = 103
$ws1.Cells.Item(245,5).Value = "LP1912"
$ws1.Cells.Item(246,1).Value = "11:43:19"
$ws1.Cells.Item(246,2).Value = "13:26"
$ws1.Cells.Item(246,3).Value = "15_ABASTO"
$ws1.Cells.Item(246,4).Value = 103
$ws1.Cells.Item(246,5).Value = "LP1912"
$ws1.Cells.Item(247,1).Value = "12:24:14"
$ws1.Cells.Item(247,2).Value = "13:27"
$ws1.Cells.Item(247,3).Value = "14_ABASTO"
$ws1.Cells.Item(247,4).Value = 63
$ws1.Cells.Item(247,5).Value = "LP1912"
$ws1.Cells.Item(248,1).Value = "11:58:46"
$ws1.Cells.Item(248,2).Value = "13:29"
$ws1.Cells.Item(248,3).Value = "17_ROMERO"
$ws1.Cells.Item(248,4).Value = 91
$ws1.Cells.Item(248,5).Value = "LP1912"
$ws1.Cells.Item(249,1).Value = "12:44:05"
$ws1.Cells.Item(249,2).Value = "13:32"
$ws1.Cells.Item(249,3).Value = "10_OLMOS"
$ws1.Cells.Item(249,4).Value = 48
$ws1.Cells.Item(249,5).Value = "LP1912"
$ws1.Cells.Item(250,1).Value = "12:57:33"
$ws1.Cells.Item(250,2).Value = "13:34"
$ws1.Cells.Item(250,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(250,4).Value = 37
$ws1.Cells.Item(250,5).Value = "LP1912"
$ws1.Cells.Item(251,1).Value = "11:43:19"
$ws1.Cells.Item(251,2).Value = "13:37"
$ws1.Cells.Item(251,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(251,4).Value = 114
$ws1.Cells.Item(251,5).Value = "LP1912"
$ws1.Cells.Item(252,1).Value = "12:24:14"
$ws1.Cells.Item(252,2).Value = "13:40"
$ws1.Cells.Item(252,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(252,4).Value = 76
$ws1.Cells.Item(252,5).Value = "LP1912"
$ws1.Cells.Item(253,1).Value = "12:44:05"
$ws1.Cells.Item(253,2).Value = "13:41"
$ws1.Cells.Item(253,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(253,4).Value = 57
$ws1.Cells.Item(253,5).Value = "LP1912"
$ws1.Cells.Item(254,1).Value = "11:58:46"
$ws1.Cells.Item(254,2).Value = "13:46"
$ws1.Cells.Item(254,3).Value = "17_ROMERO"
$ws1.Cells.Item(254,4).Value = 108
$ws1.Cells.Item(254,5).Value = "LP1912"
$ws1.Cells.Item(255,1).Value = "12:24:14"
$ws1.Cells.Item(255,2).Value = "13:47"
$ws1.Cells.Item(255,3).Value = "17_ROMERO"
$ws1.Cells.Item(255,4).Value = 83
$ws1.Cells.Item(255,5).Value = "LP1912"
$ws1.Cells.Item(256,1).Value = "11:58:46"
$ws1.Cells.Item(256,2).Value = "13:50"
$ws1.Cells.Item(256,3).Value = "215A_EL PATO"
$ws1.Cells.Item(256,4).Value = 112
$ws1.Cells.Item(256,5).Value = "LP1912"
$ws1.Cells.Item(257,1).Value = "12:57:33"
$ws1.Cells.Item(257,2).Value = "13:50"
$ws1.Cells.Item(257,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(257,4).Value = 53
$ws1.Cells.Item(257,5).Value = "LP1912"
$ws1.Cells.Item(258,1).Value = "12:24:14"
$ws1.Cells.Item(258,2).Value = "13:51"
$ws1.Cells.Item(258,3).Value = "215A_EL PATO"
$ws1.Cells.Item(258,4).Value = 87
$ws1.Cells.Item(258,5).Value = "LP1912"
$ws1.Cells.Item(259,1).Value = "11:58:46"
$ws1.Cells.Item(259,2).Value = "13:56"
$ws1.Cells.Item(259,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(259,4).Value = 118
$ws1.Cells.Item(259,5).Value = "LP1912"
$ws1.Cells.Item(260,1).Value = "11:58:46"
$ws1.Cells.Item(260,2).Value = "13:56"
$ws1.Cells.Item(260,3).Value = "225_GOMEZ"
$ws1.Cells.Item(260,4).Value = 118
$ws1.Cells.Item(260,5).Value = "LP1912"
$ws1.Cells.Item(261,1).Value = "12:24:14"
$ws1.Cells.Item(261,2).Value = "13:57"
$ws1.Cells.Item(261,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(261,4).Value = 93
$ws1.Cells.Item(261,5).Value = "LP1912"
$ws1.Cells.Item(262,1).Value = "12:44:05"
$ws1.Cells.Item(262,2).Value = "14:04"
$ws1.Cells.Item(262,3).Value = "17_ROMERO"
$ws1.Cells.Item(262,4).Value = 80
$ws1.Cells.Item(262,5).Value = "LP1912"
$ws1.Cells.Item(263,1).Value = "12:44:05"
$ws1.Cells.Item(263,2).Value = "14:05"
$ws1.Cells.Item(263,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(263,4).Value = 81
$ws1.Cells.Item(263,5).Value = "LP1912"
$ws1.Cells.Item(264,1).Value = "12:57:33"
$ws1.Cells.Item(264,2).Value = "14:16"
$ws1.Cells.Item(264,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(264,4).Value = 79
$ws1.Cells.Item(264,5).Value = "LP1912"
$ws1.Cells.Item(265,1).Value = "12:24:14"
$ws1.Cells.Item(265,2).Value = "14:17"
$ws1.Cells.Item(265,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(265,4).Value = 113
$ws1.Cells.Item(265,5).Value = "LP1912"
$ws1.Cells.Item(266,1).Value = "12:24:14"
$ws1.Cells.Item(266,2).Value = "14:20"
$ws1.Cells.Item(266,3).Value = "215C_EL PATO"
$ws1.Cells.Item(266,4).Value = 116
$ws1.Cells.Item(266,5).Value = "LP1912"
$ws1.Cells.Item(267,1).Value = "12:24:14"
$ws1.Cells.Item(267,2).Value = "14:21"
$ws1.Cells.Item(267,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(267,4).Value = 117
$ws1.Cells.Item(267,5).Value = "LP1912"
$ws1.Cells.Item(268,1).Value = "12:57:33"
$ws1.Cells.Item(268,2).Value = "14:45"
$ws1.Cells.Item(268,3).Value = "14_ABASTO"
$ws1.Cells.Item(268,4).Value = 108
$ws1.Cells.Item(268,5).Value = "LP1912"
$ws1.Cells.Item(269,1).Value = "12:57:33"
$ws1.Cells.Item(269,2).Value = "14:56"
$ws1.Cells.Item(269,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(269,4).Value = 119
$ws1.Cells.Item(269,5).Value = "LP1912"

# ---- Hoja "LP1912-215": solo cambia la hora de actualizacion ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 12:57:33"

# ---- Hoja "6203-6173": cabecera + 1 fila nueva (36 -> 37 filas) ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 12:57:33"
$ws3.Cells.Item(3,1).Value = "Total filas: 37"
$ws3.Cells.Item(42,1).Value = "12:57:33"
$ws3.Cells.Item(42,2).Value = "14:53"
$ws3.Cells.Item(42,3).Value = "215D_LA PLATA"
$ws3.Cells.Item(42,4).Value = 116
$ws3.Cells.Item(42,5).Value = "L6203"
